$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 ---
$ws.Range("A16").Value = 112178654
$ws.Range("B16").Value = 89686
$ws.Range("C16").Value = "Ovaliderad"
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 658
$ws.Range("F16").Value = "Rosenticka"
$ws.Range("G16").Value = "Rhodofomes roseus"
$ws.Range("H16").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("P16").Value = "Källåsen, Mpd"
$ws.Range("Q16").Value = 618387.2260358589
$ws.Range("R16").Value = 6904851.227267566
$ws.Range("S16").Value = 5
$ws.Range("T16").Value = "Västernorrland"
$ws.Range("U16").Value = "Sundsvall"
$ws.Range("V16").Value = "Medelpad"
$ws.Range("W16").Value = "Njurunda"
# Date-looking text must stay text (not auto-converted to a date serial) -
# enter with a leading apostrophe, then reset the style so no quote-prefix
# formatting sticks around on the cell.
$ws.Range("Y16").Value = "'2023-09-13"
$ws.Range("Y16").Style = "Normal"
$ws.Range("Z16").Value = "00:00"
$ws.Range("AA16").Value = "'2023-09-13"
$ws.Range("AA16").Style = "Normal"
$ws.Range("AB16").Value = "00:00"
$ws.Range("AD16").Value = $false
$ws.Range("AE16").Value = $false
$ws.Range("AG16").Value = $false
$ws.Range("AW16").Value = "Olle Finnström"
$ws.Range("AX16").Value = "Olle Finnström"

# --- Row 17 ---
$ws.Range("A17").Value = 112178651
$ws.Range("B17").Value = 86223
$ws.Range("C17").Value = "Ovaliderad"
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 4412
$ws.Range("F17").Value = "Äggvaxskivling"
$ws.Range("G17").Value = "Hygrophorus karstenii"
$ws.Range("H17").Value = "Sacc. & Cub."
$ws.Range("P17").Value = "Källåsen, Mpd"
$ws.Range("Q17").Value = 618387.9774688096
$ws.Range("R17").Value = 6904949.162718941
$ws.Range("S17").Value = 5
$ws.Range("T17").Value = "Västernorrland"
$ws.Range("U17").Value = "Sundsvall"
$ws.Range("V17").Value = "Medelpad"
$ws.Range("W17").Value = "Njurunda"
$ws.Range("Y17").Value = "'2023-09-13"
$ws.Range("Y17").Style = "Normal"
$ws.Range("Z17").Value = "00:00"
$ws.Range("AA17").Value = "'2023-09-13"
$ws.Range("AA17").Style = "Normal"
$ws.Range("AB17").Value = "00:00"
$ws.Range("AD17").Value = $false
$ws.Range("AE17").Value = $false
$ws.Range("AG17").Value = $false
$ws.Range("AW17").Value = "Olle Finnström"
$ws.Range("AX17").Value = "Olle Finnström"

# --- Row 18 ---
$ws.Range("A18").Value = 112178652
$ws.Range("B18").Value = 90678
$ws.Range("C18").Value = "Ovaliderad"
$ws.Range("D18").Value = "LC"
$ws.Range("E18").Value = 4366
$ws.Range("F18").Value = "Skarp dropptaggsvamp"
$ws.Range("G18").Value = "Hydnellum peckii"
$ws.Range("H18").Value = "Banker"
$ws.Range("P18").Value = "Källåsen, Mpd"
$ws.Range("Q18").Value = 618476.2382824289
$ws.Range("R18").Value = 6905001.69355389
$ws.Range("S18").Value = 5
$ws.Range("T18").Value = "Västernorrland"
$ws.Range("U18").Value = "Sundsvall"
$ws.Range("V18").Value = "Medelpad"
$ws.Range("W18").Value = "Njurunda"
$ws.Range("Y18").Value = "'2023-09-13"
$ws.Range("Y18").Style = "Normal"
$ws.Range("Z18").Value = "00:00"
$ws.Range("AA18").Value = "'2023-09-13"
$ws.Range("AA18").Style = "Normal"
$ws.Range("AB18").Value = "00:00"
$ws.Range("AD18").Value = $false
$ws.Range("AE18").Value = $false
$ws.Range("AG18").Value = $false
$ws.Range("AW18").Value = "Olle Finnström"
$ws.Range("AX18").Value = "Olle Finnström"

Write-Output "Rows 16-18 added"
